$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:B3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-12-16"
$ws.Range("B3").Value = "7.284280"
